$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = '@'
$c.Value = '42.194.84'
$c.Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  -1.08%  '
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = '@'
$c.Value = '2.240.38'
$c.Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  -1.17%  '
$ws.Cells.Item(4, 5).Value = '  +0.04%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = '@'
$c.Value = '246.26'
$c.Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -1.34%  '
$ws.Cells.Item(6, 5).Value = '  -0.46%  '
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = '@'
$c.Value = '74.65'
$c.Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  -3.02%  '
$ws.Cells.Item(8, 5).Value = '  +0.06%  '
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = '@'
$c.Value = '0.616'
$c.Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -3.47%  '
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = '@'
$c.Value = '42.30'
$c.Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +5.97%  '
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = '@'
$c.Value = '0.0941'
$c.Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -2.85%  '
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = '@'
$c.Value = '7.10'
$c.Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  -1.89%  '
$ws.Cells.Item(13, 5).Value = '  -1.68%  '
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Value = '14.47'
$c.Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -3.09%  '
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = '0.848'
$c.Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -1.50%  '
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = '@'
$c.Value = '2.240.31'
$c.Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  -1.15%  '
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = '@'
$c.Value = '42.038.17'
$c.Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  -1.24%  '
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = '@'
$c.Value = '0.0₃0984'
$c.Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  -0.60%  '
$ws.Cells.Item(19, 5).Value = '  -0.45%  '
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = '@'
$c.Value = '72.03'
$c.Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  -0.09%  '
$ws.Cells.Item(21, 2).Value = 'BitcoinCash'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = '231.27'
$c.Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -1.05%  '
$ws.Cells.Item(22, 2).Value = 'ImmutableX'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = '@'
$c.Value = '2.21'
$c.Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +3.22%  '
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = '8.83'
$c.Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +38.12%  '
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Value = '11.25'
$c.Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -0.05%  '
$ws.Cells.Item(26, 5).Value = '  -3.96%  '
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = '2.31'
$c.Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -2.44%  '
$ws.Cells.Item(28, 2).Value = 'Monero'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = '@'
$c.Value = '169.23'
$c.Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  +1.01%  '
$ws.Cells.Item(29, 2).Value = 'Toncoin'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = '@'
$c.Value = '2.09'
$c.Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -6.27%  '
$ws.Cells.Item(30, 5).Value = '  -1.26%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = '@'
$c.Value = '0.0817'
$c.Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -4.29%  '
$ws.Cells.Item(32, 5).Value = '  -2.83%  '
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = '@'
$c.Value = '30.32'
$c.Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -2.24%  '
$ws.Cells.Item(34, 5).Value = '  -1.68%  '
$ws.Cells.Item(35, 5).Value = '  +11.27%  '
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = '@'
$c.Value = '4.43'
$c.Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -2.89%  '
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = '@'
$c.Value = '0.0313'
$c.Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +2.89%  '
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = '@'
$c.Value = '13.50'
$c.Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -1.80%  '
$ws.Cells.Item(39, 5).Value = '  -3.51%  '
$ws.Cells.Item(40, 5).Value = '  -1.62%  '
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = '@'
$c.Value = '62.06'
$c.Style = 'Normal'
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = '@'
$c.Value = '0.204'
$c.Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -1.88%  '
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = '@'
$c.Value = '106.63'
$c.Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -2.33%  '
$ws.Cells.Item(44, 5).Value = '  +2.18%  '
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = '@'
$c.Value = '8.64'
$c.Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -2.32%  '
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = '@'
$c.Value = '0.996'
$c.Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -0.36%  '
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = '@'
$c.Value = '1.11'
$c.Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -3.68%  '
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = '@'
$c.Value = '4.28'
$c.Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -7.83%  '
$ws.Cells.Item(49, 5).Value = '  -0.39%  '
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = '@'
$c.Value = '2.26'
$c.Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +0.86%  '
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = '@'
$c.Value = '2.70'
$c.Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  +0.35%  '
